$d = $word.ActiveDocument
$tab = [char]9

# --- 1. "Тема курсового проекта" line: collapse the tab + 9 spaces that sat
#        between the label and the underlined title into just one extra
#        space (so "проекта " + tab + 9 spaces -> "проекта  "), retitle the
#        underlined course-project name, and collapse the trailing
#        tab + 12 spaces + tab into 8 plain spaces.

$find1 = $d.Content.Find
$ok1 = $find1.Execute("Тема курсового проекта " + $tab + "         ", $true, $false, $false, $false, $false, $true, 1, $false, "Тема курсового проекта  ", 2)
if (-not $ok1) { throw "step 1: title/tab cleanup not found" }

$find2 = $d.Content.Find
$ok2 = $find2.Execute("Мониторинг состояния системы и ядра", $true, $false, $false, $false, $false, $true, 1, $false, "Мониторинг ресурсов системы и частоты системных вызовов", 2)
if (-not $ok2) { throw "step 2: course-project title not found" }

$find3 = $d.Content.Find
$ok3 = $find3.Execute($tab + "            " + $tab, $true, $false, $false, $false, $false, $true, 1, $false, "        ", 2)
if (-not $ok3) { throw "step 3: trailing tab/space cleanup not found" }

# --- 2. Trim the slide-contents sentence: drop the trailing
#        ", результаты проведенных исследований" clause before the period.

$find4 = $d.Content.Find
$ok4 = $find4.Execute(", результаты проведенных исследований", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)
if (-not $ok4) { throw "step 4: slide-contents clause not found" }

# --- 3. Rebalance the two right-hand columns of the signatures table
#        (5773 / 2072 / 2071 -> 5773 / 2070 / 2073 twips, i.e. points
#        103.6 / 103.55 -> 103.5 / 103.65).

$t = $d.Tables.Item(2)
$t.Columns.Item(2).Width = 103.5
$t.Columns.Item(3).Width = 103.65
